# Apply updated "dSF" (column F) values to Sheet1, reflecting a repull of
# source data / recalculated means. Only column F values change; all other
# columns (including E, "dS0") remain as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0
    3  = -1
    4  = -2
    6  = -2
    9  = -3
    11 = -9
    12 = -1
    17 = 0
    25 = -1
    31 = -1
    33 = 0
    35 = -1
    38 = -1
    41 = 6
    46 = 3
    47 = -4
    51 = -3
    53 = 3
    55 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
